$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting existing rows 35..137 down to 36..138.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with this week's new data entry.
# (Same Mercado/Region/Categoria metadata as the rest of the sheet, new date & prices.)
$ws.Cells.Item(35, 1).Value = 11
$ws.Cells.Item(35, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(35, 3).Value = "Bíobío"
$ws.Cells.Item(35, 4).Value = 44742
$ws.Cells.Item(35, 4).Style = $ws.Cells.Item(36, 4).Style
$ws.Cells.Item(35, 4).NumberFormat = $ws.Cells.Item(36, 4).NumberFormat
$ws.Cells.Item(35, 5).Value = 8
$ws.Cells.Item(35, 6).Value = 100112043
$ws.Cells.Item(35, 7).Value = "Pepino ensalada"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 100
$ws.Cells.Item(35, 11).Value = 17000
$ws.Cells.Item(35, 12).Value = 18000
$ws.Cells.Item(35, 13).Value = 17500
$ws.Cells.Item(35, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(35, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(35, 16).Value = 292
$ws.Cells.Item(35, 17).Value = 60
$ws.Cells.Item(35, 18).Value = "Hortaliza"
